$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Helper/staging cell used to build the literal date-text string via a
# formula (TRIM) so that, when its *computed* value is pasted into the
# target cells as a value, Excel's "looks like a date -> convert to a
# date serial" auto-detection (which triggers on direct string literal
# assignment) is not applied. This keeps the destination cells as plain
# text matching the original inline-string formatting/style exactly.
$helper = $ws.Range("Z1")

for ($r = 3; $r -le 25; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # Column H = PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($r, 9)   # Column I = LAST UPDATE
    $helper.Formula = '=TRIM("04-Nov-2025")'
    $helper.Copy()
    $iCell.PasteSpecial(-4163)       # xlPasteValues
}

$helper.ClearContents()
$excel.CutCopyMode = $false

